$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("G1")
$ws2 = $wb.Worksheets.Item("G2")
$ws3 = $wb.Worksheets.Item("G3")

$bulletSteps = "* do this`n* do that`n* and this and that"

# --- Sheet G2: convert the plain "Steps" text to an unordered (bulleted) list ---
$ws2.Range("G4").Value = $bulletSteps
$ws2.Range("G5").Value = $bulletSteps
$ws2.Range("G6").Value = $bulletSteps
$ws2.Range("G7").Value = $bulletSteps

# --- Sheet G3: same update ---
$ws3.Range("G4").Value = $bulletSteps
$ws3.Range("G5").Value = $bulletSteps
$ws3.Range("G6").Value = $bulletSteps
$ws3.Range("G7").Value = $bulletSteps

# --- Sheet G2 / G3: widen column G a bit now that the bulleted text is there ---
$ws2.Range("G1").ColumnWidth = 18.83
$ws3.Range("G1").ColumnWidth = 18.83

# --- Sheet G2 / G3: append the same blank/wrapped rows 11-13 that G1 already has ---
$ws2.Range("G11").Value = " "
$ws2.Range("G12").Value = " "
$ws2.Range("G13").Value = " "
$ws2.Range("G12").WrapText = $true

$ws3.Range("G11").Value = " "
$ws3.Range("G12").Value = " "
$ws3.Range("G13").Value = " "
$ws3.Range("G12").WrapText = $true

# --- Update the current selection on each sheet ---
$ws1.Activate()
[void]$ws1.Range("E13").Select()

$ws2.Activate()
[void]$ws2.Range("I13").Select()

$ws3.Activate()
[void]$ws3.Range("I15").Select()

# Restore G1 as the active/visible tab (matches the original tabSelected sheet)
$ws1.Activate()
